# Generate Report for Handback
#
# The 848b503e-... file has just been handed back (in addition to the two
# files that were already handed back earlier). The localization-status
# report is regenerated: on every sheet the "848b503e" row is moved to the
# top of the data block (most-recently-handed-back first) and its
# status/target/handback fields are refreshed; the other rows keep their
# data but shift down by one row. Hyperlinks are rebuilt to match.

$wb = $excel.ActiveWorkbook

# Cell values are written with a leading apostrophe so the runtime stores
# them as literal text (matching the source t="s" shared-string cells)
# instead of auto-detecting booleans/dates/numbers (e.g. "True"/"False"
# or "2016-08-28 17:02:38") and an empty string still yields an explicit
# empty-text cell instead of clearing/removing the cell entirely.
function Set-RowValues {
    param($ws, [int]$row, [object[]]$vals)
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value2 = "'" + $vals[$i]
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-RowValues $wsOverview 2 @(
    "848b503e-2a43-4af4-bb84-f5c445444957.md",
    "e2e\848b503e-2a43-4af4-bb84-f5c445444957.md",
    ".md",
    "",
    "Handed back: in sync with en-US",
    "Handed back: in sync with en-US",
    "2016-08-28 17:02:38"
)

Set-RowValues $wsOverview 3 @(
    "ffff187ec8cc-55eb-458d-bc03-d819daf81daa.md",
    "e2e\ffff187ec8cc-55eb-458d-bc03-d819daf81daa.md",
    ".md",
    "",
    "Handed back: in sync with en-US",
    "Handed back: in sync with en-US",
    "2016-08-28 17:02:38"
)

Set-RowValues $wsOverview 4 @(
    "ffffff69b14645-b5bc-4646-8080-fef8311903ac.md",
    "e2e\ffffff69b14645-b5bc-4646-8080-fef8311903ac.md",
    ".md",
    "",
    "Handed back: in sync with en-US",
    "Handed back: in sync with en-US",
    "2016-08-28 17:02:38"
)

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0c515b3704a14e6eff4d7eefe90d5501b18c14f/e2e/848b503e-2a43-4af4-bb84-f5c445444957.md", [Type]::Missing, [Type]::Missing, "e2e\848b503e-2a43-4af4-bb84-f5c445444957.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/59ba5bf1e78592d5b6b2041913383b4c7ef225e4/e2e/ffff187ec8cc-55eb-458d-bc03-d819daf81daa.md", [Type]::Missing, [Type]::Missing, "e2e\ffff187ec8cc-55eb-458d-bc03-d819daf81daa.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0c515b3704a14e6eff4d7eefe90d5501b18c14f/e2e/ffffff69b14645-b5bc-4646-8080-fef8311903ac.md", [Type]::Missing, [Type]::Missing, "e2e\ffffff69b14645-b5bc-4646-8080-fef8311903ac.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-RowValues $wsZhCn 2 @(
    "848b503e-2a43-4af4-bb84-f5c445444957.md",
    ".md",
    "Handed back: in sync with en-US",
    "e2e",
    "ht",
    "False",
    "848b503e-2a43-4af4-bb84-f5c445444957.55a6dd8cfe27aac39e57212968d2f6922d470826.zh-cn.xlf",
    "2016-08-28 17:03:34",
    "848b503e-2a43-4af4-bb84-f5c445444957.md",
    "848b503e-2a43-4af4-bb84-f5c445444957.55a6dd8cfe27aac39e57212968d2f6922d470826.zh-cn.xlf",
    "2016-08-28 17:04:05",
    "",
    "True",
    "",
    "False",
    ""
)

Set-RowValues $wsZhCn 3 @(
    "ffff187ec8cc-55eb-458d-bc03-d819daf81daa.md",
    ".md",
    "Handed back: in sync with en-US",
    "e2e",
    "ht",
    "False",
    "cf0d3d15-07e2-4798-9d1c-82fec8a93031.ae237857e9cd159d190a9ce1b72edf3cc271d7e0.zh-cn.xlf",
    "2016-08-28 17:02:34",
    "cf0d3d15-07e2-4798-9d1c-82fec8a93031.md",
    "cf0d3d15-07e2-4798-9d1c-82fec8a93031.ae237857e9cd159d190a9ce1b72edf3cc271d7e0.zh-cn.xlf",
    "2016-08-28 17:02:51",
    "",
    "True",
    "",
    "False",
    ""
)

Set-RowValues $wsZhCn 4 @(
    "ffffff69b14645-b5bc-4646-8080-fef8311903ac.md",
    ".md",
    "Handed back: in sync with en-US",
    "e2e",
    "ht",
    "True",
    "cf0d3d15-07e2-4798-9d1c-82fec8a93031.ae237857e9cd159d190a9ce1b72edf3cc271d7e0.zh-cn.xlf",
    "2016-08-28 17:02:34",
    "cf0d3d15-07e2-4798-9d1c-82fec8a93031.md",
    "cf0d3d15-07e2-4798-9d1c-82fec8a93031.ae237857e9cd159d190a9ce1b72edf3cc271d7e0.zh-cn.xlf",
    "2016-08-28 17:02:51",
    "",
    "True",
    "",
    "False",
    ""
)

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0c515b3704a14e6eff4d7eefe90d5501b18c14f/e2e/848b503e-2a43-4af4-bb84-f5c445444957.md", [Type]::Missing, [Type]::Missing, "848b503e-2a43-4af4-bb84-f5c445444957.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0c515b3704a14e6eff4d7eefe90d5501b18c14f/e2e/848b503e-2a43-4af4-bb84-f5c445444957.md", [Type]::Missing, [Type]::Missing, "848b503e-2a43-4af4-bb84-f5c445444957.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/59ba5bf1e78592d5b6b2041913383b4c7ef225e4/e2e/ffff187ec8cc-55eb-458d-bc03-d819daf81daa.md", [Type]::Missing, [Type]::Missing, "ffff187ec8cc-55eb-458d-bc03-d819daf81daa.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/490de421cd16bbb2073620eca01285a06b80dce5/e2e/cf0d3d15-07e2-4798-9d1c-82fec8a93031.md", [Type]::Missing, [Type]::Missing, "cf0d3d15-07e2-4798-9d1c-82fec8a93031.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0c515b3704a14e6eff4d7eefe90d5501b18c14f/e2e/ffffff69b14645-b5bc-4646-8080-fef8311903ac.md", [Type]::Missing, [Type]::Missing, "ffffff69b14645-b5bc-4646-8080-fef8311903ac.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/490de421cd16bbb2073620eca01285a06b80dce5/e2e/cf0d3d15-07e2-4798-9d1c-82fec8a93031.md", [Type]::Missing, [Type]::Missing, "cf0d3d15-07e2-4798-9d1c-82fec8a93031.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-RowValues $wsDeDe 2 @(
    "848b503e-2a43-4af4-bb84-f5c445444957.md",
    ".md",
    "Handed back: in sync with en-US",
    "e2e",
    "ht",
    "False",
    "848b503e-2a43-4af4-bb84-f5c445444957.55a6dd8cfe27aac39e57212968d2f6922d470826.de-de.xlf",
    "2016-08-28 17:03:39",
    "848b503e-2a43-4af4-bb84-f5c445444957.md",
    "848b503e-2a43-4af4-bb84-f5c445444957.55a6dd8cfe27aac39e57212968d2f6922d470826.de-de.xlf",
    "2016-08-28 17:04:13",
    "",
    "True",
    "",
    "False",
    ""
)

Set-RowValues $wsDeDe 3 @(
    "ffff187ec8cc-55eb-458d-bc03-d819daf81daa.md",
    ".md",
    "Handed back: in sync with en-US",
    "e2e",
    "ht",
    "False",
    "cf0d3d15-07e2-4798-9d1c-82fec8a93031.ae237857e9cd159d190a9ce1b72edf3cc271d7e0.de-de.xlf",
    "2016-08-28 17:02:38",
    "cf0d3d15-07e2-4798-9d1c-82fec8a93031.md",
    "cf0d3d15-07e2-4798-9d1c-82fec8a93031.ae237857e9cd159d190a9ce1b72edf3cc271d7e0.de-de.xlf",
    "2016-08-28 17:02:57",
    "",
    "True",
    "",
    "False",
    ""
)

Set-RowValues $wsDeDe 4 @(
    "ffffff69b14645-b5bc-4646-8080-fef8311903ac.md",
    ".md",
    "Handed back: in sync with en-US",
    "e2e",
    "ht",
    "True",
    "cf0d3d15-07e2-4798-9d1c-82fec8a93031.ae237857e9cd159d190a9ce1b72edf3cc271d7e0.de-de.xlf",
    "2016-08-28 17:02:38",
    "cf0d3d15-07e2-4798-9d1c-82fec8a93031.md",
    "cf0d3d15-07e2-4798-9d1c-82fec8a93031.ae237857e9cd159d190a9ce1b72edf3cc271d7e0.de-de.xlf",
    "2016-08-28 17:02:57",
    "",
    "True",
    "",
    "False",
    ""
)

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0c515b3704a14e6eff4d7eefe90d5501b18c14f/e2e/848b503e-2a43-4af4-bb84-f5c445444957.md", [Type]::Missing, [Type]::Missing, "848b503e-2a43-4af4-bb84-f5c445444957.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0c515b3704a14e6eff4d7eefe90d5501b18c14f/e2e/848b503e-2a43-4af4-bb84-f5c445444957.md", [Type]::Missing, [Type]::Missing, "848b503e-2a43-4af4-bb84-f5c445444957.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/59ba5bf1e78592d5b6b2041913383b4c7ef225e4/e2e/ffff187ec8cc-55eb-458d-bc03-d819daf81daa.md", [Type]::Missing, [Type]::Missing, "ffff187ec8cc-55eb-458d-bc03-d819daf81daa.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/fad172227a2444a32d5ba10174e0e63a6fb278f6/e2e/cf0d3d15-07e2-4798-9d1c-82fec8a93031.md", [Type]::Missing, [Type]::Missing, "cf0d3d15-07e2-4798-9d1c-82fec8a93031.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0c515b3704a14e6eff4d7eefe90d5501b18c14f/e2e/ffffff69b14645-b5bc-4646-8080-fef8311903ac.md", [Type]::Missing, [Type]::Missing, "ffffff69b14645-b5bc-4646-8080-fef8311903ac.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/fad172227a2444a32d5ba10174e0e63a6fb278f6/e2e/cf0d3d15-07e2-4798-9d1c-82fec8a93031.md", [Type]::Missing, [Type]::Missing, "cf0d3d15-07e2-4798-9d1c-82fec8a93031.md") | Out-Null

$wb.Save()
